# Fix contact information missing from short resumes
#
# Inserts a new centered paragraph with the contact-info line directly
# after the "Dheeraj Chand" name paragraph (and before "PROFESSIONAL
# SUMMARY"), matching the long-resume layout.

$d = $word.ActiveDocument

function XmlEscape($s) {
    if ($null -eq $s) { return "" }
    $s = $s -replace '&', '&amp;'
    $s = $s -replace '<', '&lt;'
    $s = $s -replace '>', '&gt;'
    return $s
}

# --- Locate the name paragraph ("Dheeraj Chand") -------------------------
$searchRange = $d.Content
$found = $searchRange.Find.Execute("Dheeraj Chand", $true, $false, $false,
                                    $false, $false, $true, 1, $false,
                                    "", 0)
if (-not $found) {
    throw "Could not find the 'Dheeraj Chand' name paragraph"
}

$nameParagraph = $searchRange.Paragraphs(1)
$nameRange = $nameParagraph.Range.Duplicate

# Capture the existing direct formatting on the name run so we can
# reproduce it losslessly (rather than let a fresh run inherit it).
$isBold = $nameRange.Font.Bold
$sizePt = $nameRange.Font.Size

$nameText = $nameParagraph.Range.Text
# Strip the trailing paragraph mark (CR) / cell mark (BEL) characters.
$nameText = $nameText.TrimEnd([char]13, [char]7)
$nameText = XmlEscape $nameText

$rPrInner = ""
if ($isBold) { $rPrInner += "<w:b/>" }
if ($sizePt) { $rPrInner += '<w:sz w:val="' + [int]($sizePt * 2) + '"/>' }
$nameRPr = ""
if ($rPrInner -ne "") { $nameRPr = "<w:rPr>$rPrInner</w:rPr>" }

$contactText = XmlEscape "202.550.7110 | dheeraj.chand@gmail.com | https://www.dheerajchand.com | https://www.linkedin.com/in/dheerajchand/ | Austin, TX"

# Rebuild the name paragraph exactly as it was, followed by the new
# centered contact-info paragraph (no direct run formatting), and hand
# the whole fragment to InsertXML over the *original* paragraph's range.
# Doing the replace this way (old-paragraph-XML + new-paragraph-XML)
# keeps every other paragraph in the document completely untouched,
# whereas inserting at a collapsed range right on a paragraph boundary
# can corrupt the following paragraph.
$payload = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
  '<pkg:xmlData>' +
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:body>' +
  "<w:p><w:pPr><w:jc w:val=`"center`"/></w:pPr><w:r>$nameRPr<w:t>$nameText</w:t></w:r></w:p>" +
  "<w:p><w:pPr><w:jc w:val=`"center`"/></w:pPr><w:r><w:t>$contactText</w:t></w:r></w:p>" +
  '</w:body></w:document>' +
  '</pkg:xmlData></pkg:part></pkg:package>'

$nameRange.InsertXML($payload)
